# Backend Added For Marksheet 1
#
# The source marksheet's header row used "kendra_code" (column M) and
# "working_education" (column X). The backend rename swaps these for the
# actual field names used by the data pipeline:
#   M1: kendra_code          -> examination_center_code
#   X1: working_education    -> work_education
#
# All other header cells (school_name, district, block, school_dice_code,
# scholar_no, roll_no, student_name, father_name, mother_name, dob,
# student_class, swayam_pathi, marksheet_id, total_grade, examination_date,
# hindi, english, maths, sanskrit, science, social_science,
# environmental_studies, physical, arts) stay exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "examination_center_code"
$ws.Range("X1").Value = "work_education"

# Leave the cursor/selection where the author last left it in the sheet.
$ws.Range("R5").Select()
